$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "created"
$ws.Range("B1").Value = "was"
$ws.Range("D1").Value = "with"
$ws.Range("E1").Value = "love"

$ws.Range("A2").Value = "was"
$ws.Range("A3").Value = "created"
$ws.Range("A4").Value = "with"
$ws.Range("A5").Value = "love"

$ws.Range("C1").Select()
